# Update Data Analysis file. Update data from experiment.
#
# Adds a new "WMM" worksheet (a copy of the existing "Sheet2" experiment
# template) populated with newly measured displacement data, and makes
# it the active sheet.

$wb = $excel.ActiveWorkbook

# "WMM" starts life as a duplicate of Sheet2's layout/styles (same table
# of trial/replicate/displacement data), inserted right after Sheet2.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Copy($null, $ws2)
$wmm = $wb.Worksheets.Item($wb.Worksheets.Count)
$wmm.Name = "WMM"

# This sheet records pulse width in microseconds rather than milliseconds.
$wmm.Range("A3").Value = "t (us):"

# Newly measured displacement values (cm) for the Flexor/Ulnar trials.
$wmm.Range("D9").Value  = 1.7
$wmm.Range("E9").Value  = 1
$wmm.Range("F9").Value  = 0.4

$wmm.Range("D10").Value = 1.4
$wmm.Range("E10").Value = 0.5
$wmm.Range("F10").Value = 0.8

$wmm.Range("D11").Value = 1.4
$wmm.Range("E11").Value = 1
$wmm.Range("F11").Value = 0.4

$wmm.Range("D12").Value = 0.9
$wmm.Range("E12").Value = 0.5
$wmm.Range("F12").Value = 0.4

$wmm.Range("D13").Value = 1
$wmm.Range("E13").Value = 1
$wmm.Range("F13").Value = 0.7

$wmm.Range("D14").Value = 1.25
$wmm.Range("E14").Value = 0.5
$wmm.Range("F14").Value = 0.4

$wmm.Range("D15").Value = 3.3
$wmm.Range("E15").Value = 0.6
$wmm.Range("F15").Value = 0.75

$wmm.Range("D16").Value = 4.5
$wmm.Range("E16").Value = 1.1
$wmm.Range("F16").Value = 1.9

$wmm.Range("D17").Value = 4.8
$wmm.Range("E17").Value = 1.2
$wmm.Range("F17").Value = 1.3

$wmm.Range("E18").Value = 0.5
$wmm.Range("F18").Value = 1.2

$wmm.Range("E19").Value = 0.4
$wmm.Range("F19").Value = 1.3

$wmm.Range("E20").Value = 0.4
$wmm.Range("F20").Value = 1.2

# Land the selection/active tab on the new sheet, as the author left it.
$wmm.Range("A4").Select()
$wmm.Activate()
